$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new row (14) for 아이유 (IU) to the singers table.
$ws.Range("A14").Value = "13"
$ws.Range("B14").Value = "아이유"
$ws.Range("C14").Value = "1993년 05월 16일"
$ws.Range("E14").Value = "EDAM엔터테인먼트"
$ws.Range("F14").Value = "대한민국"

# Update the active selection to mirror the saved view state.
$ws.Range("F17").Select()
